# The sheet has 7 columns: A=code, B=name, C=status,
# D=codeforiati:group-name, E=codeforiati:category-name,
# F=codeforiati:group-code, G=codeforiati:category-code.
#
# The edit reorders these last four columns (header + every data row) to:
# D=codeforiati:category-code, E=codeforiati:group-code,
# F=codeforiati:group-name,    G=codeforiati:category-name
#
# i.e. (D,E,F,G) -> (G,F,D,E). The underlying group/category code & name
# values for every row stay exactly the same, they just live under
# differently-ordered headers afterwards.
#
# We implement the rearrangement with Copy / PasteSpecial (values only)
# through a scratch staging area, rather than Range.Value assignment,
# so that cells keep their original text/shared-string representation
# (this matters because several of the codes, e.g. "110", "111", look
# numeric and would otherwise be silently re-typed as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$srcRange  = "D1:G" + $lastRow
$stageD    = "I1:I" + $lastRow
$stageE    = "J1:J" + $lastRow
$stageF    = "K1:K" + $lastRow
$stageG    = "L1:L" + $lastRow
$stageAll  = "I1:L" + $lastRow

# 1) Stash the current D,E,F,G columns into the scratch columns I,J,K,L.
$ws.Range($srcRange).Copy()
$ws.Range("I1").PasteSpecial(-4163)

# 2) new D = old G (stage L) ; new E = old F (stage K)
#    new F = old D (stage I) ; new G = old E (stage J)
$ws.Range($stageG).Copy()
$ws.Range("D1").PasteSpecial(-4163)

$ws.Range($stageF).Copy()
$ws.Range("E1").PasteSpecial(-4163)

$ws.Range($stageD).Copy()
$ws.Range("F1").PasteSpecial(-4163)

$ws.Range($stageE).Copy()
$ws.Range("G1").PasteSpecial(-4163)

# 3) Clean up the scratch area.
$ws.Range($stageAll).ClearContents()
